# Auto-generated edit script: updates numeric cells in columns H-N
# (currentAveragePrice*, LevePrice*, LeveProfit*) across multiple sheets
# to match the target workbook state.

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 31
$ws.Range("H31").Value = 267.33334
$ws.Range("I31").Value = 267.33334
$ws.Range("K31").Value = 802.0000200000001
$ws.Range("M31").Value = -572.0000200000001
# Row 40
$ws.Range("H40").Value = 2025.5
$ws.Range("I40").Value = 888.375
$ws.Range("J40").Value = 3162.625
$ws.Range("K40").Value = 888.375
$ws.Range("L40").Value = 3162.625
$ws.Range("M40").Value = -713.375
$ws.Range("N40").Value = -3512.625
# Row 48
$ws.Range("H48").Value = 4999
$ws.Range("J48").Value = 4999
$ws.Range("L48").Value = 14997
$ws.Range("N48").Value = -15581
# Row 56
$ws.Range("H56").Value = 4999
$ws.Range("J56").Value = 4999
$ws.Range("L56").Value = 14997
$ws.Range("N56").Value = -16065
# Row 70
$ws.Range("H70").Value = 1000
$ws.Range("I70").Value = 0
$ws.Range("J70").Value = 1000
$ws.Range("K70").Value = 0
$ws.Range("L70").Value = 3000
$ws.Range("M70").ClearContents()
$ws.Range("N70").Value = -3540
# Row 73
$ws.Range("H73").Value = 1000
$ws.Range("I73").Value = 0
$ws.Range("J73").Value = 1000
$ws.Range("K73").Value = 0
$ws.Range("L73").Value = 3000
$ws.Range("M73").ClearContents()
$ws.Range("N73").Value = -4872
# Row 111
$ws.Range("H111").Value = 2995
$ws.Range("I111").Value = 2995
$ws.Range("K111").Value = 8985
$ws.Range("M111").Value = -5918
# Row 113
$ws.Range("H113").Value = 4315
$ws.Range("I113").Value = 2670.8
$ws.Range("J113").Value = 4771.722
$ws.Range("K113").Value = 2670.8
$ws.Range("L113").Value = 4771.722
$ws.Range("M113").Value = 583.1999999999998
$ws.Range("N113").Value = -11279.722
# Row 115
$ws.Range("H115").Value = 4334
$ws.Range("I115").Value = 2096
$ws.Range("K115").Value = 6288
$ws.Range("M115").Value = -4721
# Row 116
$ws.Range("H116").Value = 9799
$ws.Range("I116").Value = 9749.166999999999
$ws.Range("K116").Value = 9749.166999999999
$ws.Range("M116").Value = -6307.166999999999
# Row 132
$ws.Range("H132").Value = 3033.4614
$ws.Range("I132").Value = 1079.9697
$ws.Range("K132").Value = 3239.9091
$ws.Range("M132").Value = -709.9091000000003
# Row 137
$ws.Range("H137").Value = 64637.562
$ws.Range("I137").Value = 2667.4285
$ws.Range("J137").Value = 112836.555
$ws.Range("K137").Value = 8002.2855
$ws.Range("L137").Value = 338509.665
$ws.Range("M137").Value = -5452.2855
$ws.Range("N137").Value = -343609.665
# Row 138
$ws.Range("H138").Value = 2739.6155
$ws.Range("I138").Value = 1415.25
$ws.Range("J138").Value = 3328.2222
$ws.Range("K138").Value = 4245.75
$ws.Range("L138").Value = 9984.6666
$ws.Range("M138").Value = 894.25
$ws.Range("N138").Value = -20264.6666

$ws = $wb.Worksheets.Item("ARM")
# Row 61
$ws.Range("H61").Value = 4366.3335
$ws.Range("I61").Value = 4366.3335
$ws.Range("K61").Value = 4366.3335
$ws.Range("M61").Value = -4154.3335
# Row 107
$ws.Range("H107").Value = 44999
$ws.Range("J107").Value = 44999
$ws.Range("L107").Value = 44999
$ws.Range("N107").Value = -52679
# Row 132
$ws.Range("H132").Value = 4335.6665
$ws.Range("I132").Value = 0
$ws.Range("J132").Value = 4335.6665
$ws.Range("K132").Value = 0
$ws.Range("L132").Value = 13006.9995
$ws.Range("M132").ClearContents()
$ws.Range("N132").Value = -18066.9995
# Row 136
$ws.Range("H136").Value = 4366.3335
$ws.Range("I136").Value = 4366.3335
$ws.Range("K136").Value = 13099.0005
$ws.Range("M136").Value = -10549.0005

$ws = $wb.Worksheets.Item("BSM")
# Row 20
$ws.Range("H20").Value = 14240.615
$ws.Range("J20").Value = 15384.571
$ws.Range("L20").Value = 15384.571
$ws.Range("N20").Value = -15878.571
# Row 87
$ws.Range("H87").Value = 0
$ws.Range("J87").Value = 0
$ws.Range("L87").Value = 0
$ws.Range("N87").ClearContents()
# Row 90
$ws.Range("H90").Value = 0
$ws.Range("J90").Value = 0
$ws.Range("L90").Value = 0
$ws.Range("N90").ClearContents()
# Row 92
$ws.Range("H92").Value = 0
$ws.Range("J92").Value = 0
$ws.Range("L92").Value = 0
$ws.Range("N92").ClearContents()
# Row 94
$ws.Range("H94").Value = 1349
$ws.Range("I94").Value = 1049.6666
$ws.Range("K94").Value = 1049.6666
$ws.Range("M94").Value = -598.6666
# Row 95
$ws.Range("H95").Value = 0
$ws.Range("J95").Value = 0
$ws.Range("L95").Value = 0
$ws.Range("N95").ClearContents()
# Row 105
$ws.Range("H105").Value = 5267495.5
$ws.Range("I105").Value = 6671653.5
$ws.Range("K105").Value = 6671653.5
$ws.Range("M105").Value = -6669906.5
# Row 107
$ws.Range("H107").Value = 14406.863
$ws.Range("I107").Value = 3525.389
$ws.Range("K107").Value = 3525.389
$ws.Range("M107").Value = -1605.389
# Row 109
$ws.Range("H109").Value = 79789.336
$ws.Range("J109").Value = 79789.336
$ws.Range("L109").Value = 79789.336
$ws.Range("N109").Value = -82563.336
# Row 110
$ws.Range("H110").Value = 0
$ws.Range("J110").Value = 0
$ws.Range("L110").Value = 0
$ws.Range("N110").ClearContents()
# Row 111
$ws.Range("H111").Value = 0
$ws.Range("J111").Value = 0
$ws.Range("L111").Value = 0
$ws.Range("N111").ClearContents()
# Row 134
$ws.Range("H134").Value = 2107.2666
$ws.Range("I134").Value = 1900.6428
$ws.Range("K134").Value = 5701.928400000001
$ws.Range("M134").Value = -3166.928400000001

$ws = $wb.Worksheets.Item("CRP")
# Row 31
$ws.Range("H31").Value = 4473.75
$ws.Range("I31").Value = 4473.75
$ws.Range("J31").Value = 0
$ws.Range("K31").Value = 4473.75
$ws.Range("L31").Value = 0
$ws.Range("M31").Value = -4178.75
$ws.Range("N31").ClearContents()
# Row 34
$ws.Range("H34").Value = 4473.75
$ws.Range("I34").Value = 4473.75
$ws.Range("J34").Value = 0
$ws.Range("K34").Value = 4473.75
$ws.Range("L34").Value = 0
$ws.Range("M34").Value = -4271.75
$ws.Range("N34").ClearContents()
# Row 86
$ws.Range("H86").Value = 8158.154
$ws.Range("I86").Value = 7707.5557
$ws.Range("K86").Value = 7707.5557
$ws.Range("M86").Value = -6584.5557
# Row 89
$ws.Range("H89").Value = 8158.154
$ws.Range("I89").Value = 7707.5557
$ws.Range("K89").Value = 38537.7785
$ws.Range("M89").Value = -32921.7785
# Row 105
$ws.Range("H105").Value = 3579.8572
$ws.Range("I105").Value = 4036.8333
$ws.Range("J105").Value = 3096
$ws.Range("K105").Value = 4036.8333
$ws.Range("L105").Value = 3096
$ws.Range("M105").Value = -2289.8333
$ws.Range("N105").Value = -6590
# Row 107
$ws.Range("H107").Value = 770.7895
$ws.Range("I107").Value = 463.14285
$ws.Range("K107").Value = 463.14285
$ws.Range("M107").Value = 1456.85715
# Row 132
$ws.Range("H132").Value = 2347.1428
$ws.Range("I132").Value = 2489.4
$ws.Range("J132").Value = 1991.5
$ws.Range("K132").Value = 7468.200000000001
$ws.Range("L132").Value = 5974.5
$ws.Range("M132").Value = -4938.200000000001
$ws.Range("N132").Value = -11034.5
# Row 134
$ws.Range("H134").Value = 3300
$ws.Range("I134").Value = 3246.1538
$ws.Range("J134").Value = 4000
$ws.Range("K134").Value = 9738.4614
$ws.Range("L134").Value = 12000
$ws.Range("M134").Value = -7203.4614
$ws.Range("N134").Value = -17070

$ws = $wb.Worksheets.Item("CUL")
# Row 4
$ws.Range("H4").Value = 7867852.5
$ws.Range("I4").Value = 1214850.9
$ws.Range("K4").Value = 3644552.7
$ws.Range("M4").Value = -3644440.7

$ws = $wb.Worksheets.Item("GSM")
# Row 102
$ws.Range("H102").Value = 2751.6875
$ws.Range("I102").Value = 2957
$ws.Range("K102").Value = 2957
$ws.Range("M102").Value = -1335
# Row 122
$ws.Range("H122").Value = 1151.4286
$ws.Range("I122").Value = 1151.4286
$ws.Range("K122").Value = 3454.2858
$ws.Range("M122").Value = -1004.2858
# Row 126
$ws.Range("H126").Value = 2200
$ws.Range("I126").Value = 900
$ws.Range("K126").Value = 2700
$ws.Range("M126").Value = -230
# Row 132
$ws.Range("H132").Value = 800
$ws.Range("I132").Value = 800
$ws.Range("J132").Value = 0
$ws.Range("K132").Value = 2400
$ws.Range("L132").Value = 0
$ws.Range("M132").Value = 130
$ws.Range("N132").ClearContents()

$ws = $wb.Worksheets.Item("LTW")
# Row 22
$ws.Range("H22").Value = 1489.125
$ws.Range("J22").Value = 1166.6666
$ws.Range("L22").Value = 1166.6666
$ws.Range("N22").Value = -1756.6666
# Row 27
$ws.Range("H27").Value = 1489.125
$ws.Range("J27").Value = 1166.6666
$ws.Range("L27").Value = 1166.6666
$ws.Range("N27").Value = -1380.6666
# Row 46
$ws.Range("H46").Value = 54908.375
$ws.Range("I46").Value = 86193.8
$ws.Range("K46").Value = 86193.8
$ws.Range("M46").Value = -86005.8
# Row 132
$ws.Range("H132").Value = 4876.7026
$ws.Range("I132").Value = 4966.9688
$ws.Range("K132").Value = 14900.9064
$ws.Range("M132").Value = -12370.9064

$ws = $wb.Worksheets.Item("WVR")
# Row 38
$ws.Range("H38").Value = 19530.5
$ws.Range("J38").Value = 19530.5
$ws.Range("L38").Value = 19530.5
$ws.Range("N38").Value = -20476.5
# Row 49
$ws.Range("H49").Value = 33495
$ws.Range("J49").Value = 33495
$ws.Range("L49").Value = 33495
$ws.Range("N49").Value = -33955
# Row 132
$ws.Range("H132").Value = 5722.577
$ws.Range("I132").Value = 7551.9414
$ws.Range("J132").Value = 2267.111
$ws.Range("K132").Value = 22655.8242
$ws.Range("L132").Value = 6801.333
$ws.Range("M132").Value = -20125.8242
$ws.Range("N132").Value = -11861.333
# Row 136
$ws.Range("H136").Value = 2714.7222
$ws.Range("I136").Value = 2286.1765
$ws.Range("K136").Value = 6858.529500000001
$ws.Range("M136").Value = -4308.529500000001

